$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 updates
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 2.2
$ws.Range("J8").Value = 3.6
$ws.Range("L8").Value = 2.72
$ws.Range("W8").Value = 8.75
$ws.Range("Y8").Value = 11.25
$ws.Range("AA8").Value = 29
$ws.Range("AB8").Value = 37
$ws.Range("AD8").Value = 6.2
$ws.Range("AE8").Value = 14
$ws.Range("AH8").Value = 7.5
$ws.Range("AI8").Value = 10.75
$ws.Range("AJ8").Value = 8.75
$ws.Range("AK8").Value = 22
$ws.Range("AL8").Value = 18
$ws.Range("AM8").Value = 28
$ws.Range("AO8").Value = 17
$ws.Range("AQ8").Value = 80
$ws.Range("AT8").Value = 2.6
$ws.Range("AU8").Value = 6.8
$ws.Range("AV8").Value = 55
$ws.Range("AW8").Value = 4.1
$ws.Range("AX8").Value = 11
$ws.Range("AY8").Value = 18
$ws.Range("AZ8").Value = 40
$ws.Range("BA8").Value = 70
$ws.Range("BB8").Value = 200

# Row 9 updates
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 3
$ws.Range("Q9").Value = 2.15
$ws.Range("R9").Value = 1.67
